$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before I, so the existing "costs"/"Costs" column (I) moves to J,
# and the new "digestion"/"Digestion" column takes the I position.
$ws.Range("I1").EntireColumn.Insert()

# Header labels for new column I (digestion / Digestion)
$ws.Range("I1").Value = "digestion"
$ws.Range("I2").Value = "Digestion"

# Fill in the numeric performance values for E3:J8 (replacing the #N/A placeholders)
$data = @(
    @(0, 1, 1, -1, -1, 0),
    @(0, 0, 0, 0, 0, 0),
    @(0, 1, 0, -1, -1, 0),
    @(0, 0, -1, -1, 0, 0),
    @(0, 1, 0, -1, 0, 1),
    @(1, 0, 1, 1, 1, -1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 3 + $i
    $rowData = $data[$i]
    [object[,]]$arr = New-Object 'object[,]' 1,6
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $arr[0, $j] = $rowData[$j]
    }
    $ws.Range("E$row`:J$row").Value = $arr
}

# Update the active selection to match the saved view state
$ws.Range("I7").Select()
